# Update "想去人数" (interested-people count) figures in column F
# for the sheets "展览", "演出", and "全部类型", reflecting a refreshed
# scrape of source data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1136
$ws1.Range("F6").Value = 135
$ws1.Range("F7").Value = 59
$ws1.Range("F8").Value = 60
$ws1.Range("F10").Value = 5125
$ws1.Range("F11").Value = 4748
$ws1.Range("F13").Value = 34
$ws1.Range("F16").Value = 182

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 74

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1136
$ws4.Range("F6").Value = 135
$ws4.Range("F7").Value = 59
$ws4.Range("F8").Value = 60
$ws4.Range("F10").Value = 5125
$ws4.Range("F11").Value = 4748
$ws4.Range("F13").Value = 34
$ws4.Range("F16").Value = 182
$ws4.Range("F17").Value = 74
